# 2017-01-31 update: energy.gov - chunk 7
# Table 6.1.B - roll the "Net Summer Capacity for Estimated Distributed
# Solar Photovoltaic Capacity" table forward one month: add the
# "November" 2016 data row (just above the footnote row) and refresh the
# subtitle to match the new date range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the subtitle ("2014 - October 2016" -> "2014 - November 2016") ---
$ws.Cells.Item(2, 1).Value2 = "2014 - November 2016"

# --- Insert a new row just above the footnote row (currently row 44) ---
# This pushes the footnote row (and its A44:E44 merge) down to row 45.
$ws.Rows.Item(44).Insert()

# New row should look like the other monthly data rows, so clone the
# formatting from the row directly above it (October 2016, row 43).
$ws.Range("A43:E43").Copy()
$ws.Range("A44:E44").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Populate the November 2016 figures ---
$ws.Cells.Item(44, 1).Value2 = "November"
$ws.Cells.Item(44, 2).Value2 = 7241.5
$ws.Cells.Item(44, 3).Value2 = 4572.8999999999996
$ws.Cells.Item(44, 4).Value2 = 1049.0999999999999
$ws.Cells.Item(44, 5).Value2 = 12863.5
